$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns D and E so numeric-looking strings
# (e.g. "56.307.01", "0.430") are preserved as text, not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '56.307.01'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '2.996.64'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '506.56'
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("D6").Value = '137.84'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.430'
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").Value = '7.12'
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").Value = '0.366'
$ws.Range("E11").Value = '  +2.66%  '
$ws.Range("D12").Value = '3.508.28'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").Value = '25.55'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '56.286.09'
$ws.Range("E16").Value = '  -1.47%  '
$ws.Range("D17").Value = '2.992.15'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '5.98'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").Value = '12.94'
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = '8.05'
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("D21").Value = '331.76'
$ws.Range("E21").Value = '  +3.56%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").Value = '0.496'
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = '65.03'
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("D25").Value = '3.118.77'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = '0.165'
$ws.Range("E26").Value = '  +1.40%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0940'
$ws.Range("E28").Value = '  +5.65%  '
$ws.Range("D29").Value = '6.34'
$ws.Range("E29").Value = '  -4.14%  '
$ws.Range("D30").Value = '6.90'
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("D31").Value = '1.78'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").Value = '20.28'
$ws.Range("E32").Value = '  +0.70%  '
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").Value = '153.87'
$ws.Range("E34").Value = '  -0.36%  '
$ws.Range("D35").Value = '4.47'
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("E36").Value = '  +1.16%  '
$ws.Range("D37").Value = '26.23'
$ws.Range("E37").Value = '  +7.47%  '
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("E39").Value = '  -0.36%  '
$ws.Range("D40").Value = '3.034.23'
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("D41").Value = '36.69'
$ws.Range("E41").Value = '  -2.92%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '3.79'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").Value = '0.651'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").Value = '2.180.15'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").Value = '1.34'
$ws.Range("E46").Value = '  -2.61%  '
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").Value = '0.927'
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").Value = '19.41'
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("D51").Value = '0.0852'
$ws.Range("E51").Value = '  -1.90%  '
